$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 3

$ws.Range("G2:G6").Value = 0.5587383333333333
$ws.Range("H2:H6").Value = 1.676215

$ws.Range("K2").Value = 3
$ws.Range("K3").Value = 3
$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = 2

$ws.Range("M2").Value = 62.07563766666667
$ws.Range("N2").Value = 186.226913
$ws.Range("O2").Value = 0.4556169394345516
$ws.Range("P2").Value = 0.53808222397892
$ws.Range("Q2").Value = 34.68403833047722
$ws.Range("R2").Value = 312.156344974295
$ws.Range("S2").Value = 0.4556169394345516
$ws.Range("T2").Value = 0.53808222397892

$ws.Range("M3").Value = 8.371752000000001
$ws.Range("N3").Value = 25.115256
$ws.Range("O3").Value = 0.06144619962548196
$ws.Range("P3").Value = 0.07256777544435759
$ws.Range("Q3").Value = 4.677618759560001
$ws.Range("R3").Value = 42.09856883604
$ws.Range("S3").Value = 0.06144619962548196
$ws.Range("T3").Value = 0.07256777544435759

$ws.Range("M4").Value = 0.4371553333333333
$ws.Range("N4").Value = 1.311466
$ws.Range("O4").Value = 0.003208591687778628
$ws.Range("P4").Value = 0.003789337054374833
$ws.Range("Q4").Value = 0.2442554423544445
$ws.Range("R4").Value = 2.19829898119
$ws.Range("S4").Value = 0.003208591687778628
$ws.Range("T4").Value = 0.003789337054374833

$ws.Range("M5").Value = 2.718766
$ws.Range("N5").Value = 8.156298
$ws.Range("O5").Value = 0.01995494352567695
$ws.Range("P5").Value = 0.0235667277976885
$ws.Range("Q5").Value = 1.519078783563333
$ws.Range("R5").Value = 13.67170905207
$ws.Range("S5").Value = 0.01995494352567695
$ws.Range("T5").Value = 0.0235667277976885

$ws.Range("M6").Value = 62.6419255
$ws.Range("N6").Value = 125.283851
$ws.Range("O6").Value = 0.4597733257265108
$ws.Range("P6").Value = 0.3619939357246589
$ws.Range("Q6").Value = 35.00044505066084
$ws.Range("R6").Value = 210.002670303965
$ws.Range("S6").Value = 0.4597733257265108
$ws.Range("T6").Value = 0.3619939357246589
